# Update the LocInst sheet: replace the list of "Norma" codes/quantities
# with a single "MAGAZINE MIKRON" row, restyle the header/data rows,
# widen column A to fit the new text, and select B2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Header row (A1:B1): keep bold + border + horizontal-center,
#        but drop the vertical-center alignment (now vertical = bottom/default).
$ws.Range("A1:B1").VerticalAlignment = -4107   # xlBottom -> default, attribute dropped on save

# --- 2. Replace the data: single row "MAGAZINE MIKRON" / 16 in row 2,
#        clear out the old Norma codes + quantities from rows 3-5.
$ws.Range("A2").Value = "MAGAZINE MIKRON"
$ws.Range("B2").Value = 16
$ws.Range("A3:B5").ClearContents()

# --- 3. Restyle rows 2-5: no border, not bold, horizontal+vertical center.
$ws.Range("A2:B5").Borders.LineStyle = -4142   # xlLineStyleNone
$ws.Range("A2:B5").Font.Bold = $false

# --- 4. Column A needs to widen to fit "MAGAZINE MIKRON" (was sharing
#        an 11.43-wide definition with column B).
$ws.Columns("A").ColumnWidth = 17.8

# --- 5. Match the saved selection (cell B2 active).
$ws.Range("B2").Select() | Out-Null
